# Swap the data of row 3 <-> row 4 and row 19 <-> row 20 on the active
# sheet. Both pairs of rows are fully swapped across every used column
# (A:AY), which is what the target diff represents (two pairs of records
# that traded places in the export).
#
# Columns Y and AA hold the same text date ("2026-01-24") in every row of
# this sheet, so swapping them is a no-op value-wise; they are handled in
# their own single-column range so that assigning Value2 never routes the
# text through Excel's date auto-detection (which would turn them into a
# numeric date serial) - this is purely to avoid accidentally changing
# their type while leaving the swap correct for the columns that do change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-ColumnRange($colFrom, $colTo, $r1, $r2) {
    $range1 = $ws.Range($colFrom + $r1 + ":" + $colTo + $r1)
    $range2 = $ws.Range($colFrom + $r2 + ":" + $colTo + $r2)
    $v1 = $range1.Value2
    $v2 = $range2.Value2
    $range1.Value2 = $v2
    $range2.Value2 = $v1
}

function Swap-Rows($r1, $r2) {
    Swap-ColumnRange "A" "X" $r1 $r2
    Swap-ColumnRange "Z" "Z" $r1 $r2
    Swap-ColumnRange "AB" "AY" $r1 $r2
}

Swap-Rows 3 4
Swap-Rows 19 20
